# TC04_Trials_Filter_Race-NotReported.xlsx
# "all single filter scripts in CTDC"
# Adds a third row (FilesTab) to the "startup" sheet and refreshes the
# CasesTab / StatQuery Cypher text in row 2 to the newer single-filter form.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- updated CasesTab query (row 2, column B) ---------------------------
$casesTabQuery = @'
MATCH (c:case)
 MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)
 MATCH (f:file)-[*]->(c)
   WHERE c.race = "NOT_REPORTED"
RETURN DISTINCT
    c.case_id AS `Case ID`,
     ct.clinical_trial_designation AS `Trial Code`,
     a.arm_id AS Arm,
      a.arm_drug AS `Arm Treatment`,
c.disease AS Diagnosis,
  c.gender AS Gender,
    c.race AS Race,
    c.ethnicity AS Ethnicity
'@

# --- new FilesTab query (row 3, column B) --------------------------------
$filesTabQuery = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
OPTIONAL MATCH (f)-->(parent)
WITH f,a,ct,c,parent
        WHERE c.race = "NOT_REPORTED"
WITH
    f, parent, c, a, ct,
    ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
    toInteger(floor(log(f.file_size)/log(1024))) as i,
    2 as precision
WITH
    f, parent, c, a, ct,
    f.file_size /(1024^i) AS value,
    10^precision AS factor,
    units[i] as unit
WITH
    f, parent, c, a, ct, unit,
    round(factor * value)/factor AS size
RETURN DISTINCT
    f.file_name AS `File Name`,
    head(labels(parent)) as Association,
    f.file_description AS Description,
    f.file_format AS `File Format`,
    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    ct.clinical_trial_designation AS `Trial Code`,
    a.arm_id AS Arm,
    c.case_id AS `Case ID`
'@

# --- updated StatQuery query (shared by row 2 and row 3, column C) ------
$statQuery = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
      WHERE c.race = "NOT_REPORTED"
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files
'@

# Write order matters: it determines the order new entries land in
# xl/sharedStrings.xml once the now-unreferenced old query text is
# garbage-collected on save (FilesTab, then CasesTab query, then FilesTab
# query, then StatQuery - matching the target shared-string table).
$ws.Range("A3").Value = "FilesTab"
$ws.Range("B2").Value = $casesTabQuery
$ws.Range("B3").Value = $filesTabQuery
$ws.Range("C2").Value = $statQuery
$ws.Range("C3").Value = $statQuery

$ws.Range("D3").Value = "TC04_Trials_Filter_Race-NotReported_Neo4jData.xlsx"
$ws.Range("E3").Value = "TC04_Trials_Filter_Race-NotReported_WebData.xlsx"

$ws.Range("B2").WrapText = $true
$ws.Range("C2").WrapText = $true
$ws.Range("B3").WrapText = $true
$ws.Range("C3").WrapText = $true

# --- row heights ----------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 195
$ws.Rows.Item(3).RowHeight = 409.5

# --- view: scroll/zoom to show the new row, matching the saved state -----
$excel.ActiveWindow.Zoom = 70
$ws.Range("B3").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 2
